$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 126, pushing existing rows 126-158 down to 127-159
$ws.Rows("126").Insert()

# Fill in the new row 126 with the new weekly data point
$ws.Cells.Item(126, 1).Value = 4
$ws.Cells.Item(126, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(126, 3).Value = "Los Lagos"
$ws.Cells.Item(126, 4).Value = 44511
$ws.Cells.Item(126, 4).NumberFormat = $ws.Cells.Item(127, 4).NumberFormat
$ws.Cells.Item(126, 5).Value = 10
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100108
$ws.Cells.Item(126, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(126, 9).Value = 100108005
$ws.Cells.Item(126, 10).Value = "Piña"
$ws.Cells.Item(126, 11).Value = "Caramelo"
$ws.Cells.Item(126, 12).Value = "Segunda"
$ws.Cells.Item(126, 13).Value = 90
$ws.Cells.Item(126, 14).Value = 22000
$ws.Cells.Item(126, 15).Value = 23000
$ws.Cells.Item(126, 16).Value = 22500
$ws.Cells.Item(126, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(126, 18).Value = "Ecuador"
$ws.Cells.Item(126, 19).Value = 1607
$ws.Cells.Item(126, 20).Value = 14
